$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Format B1: bold font, thin border all around, centered horizontally, top vertically
$cell = $ws.Range("B1")
$cell.Font.Bold = $true
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4160
$cell.Borders.LineStyle = 1
$cell.Borders.Weight = 2

# Apply the exact same style to A2 by copying formats, so the same
# cell style index gets reused instead of creating a duplicate one.
$cell.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
